# "fixes undefined for chair"
#
# The "Sets" sheet had D Jasper / Rio Tinto entered as Chair/Institution
# on the "Closure Planning and Stakeholders" row (row 8) - but that data
# actually belongs to the untitled Set 8 row (row 9, the "Break" row)
# which was missing its Title. Clear the wrongly-placed Chair/Institution
# values and give row 9 its "Break" title instead. Also backfills the
# missing "End Time" for the last Sessions row, and updates the saved
# window/selection state (active sheet + selected cell) on both sheets.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sessions")
$ws2 = $wb.Worksheets.Item("Sets")

# --- Sessions sheet -------------------------------------------------
# Row 14 ("Close of day") was missing an End Time; fill it in as one
# minute after the Start Time, matching the date/time formatting
# already used by the Start Time cell (C14).
$ws1.Range("C14").Copy()
$ws1.Range("D14").PasteSpecial(-4122)
$ws1.Range("D14").Value = 43713.4798611111
$ws1.Rows.Item(14).RowHeight = 13.8

# --- Sets sheet -------------------------------------------------------
# Clear the erroneous Chair / Institution on the
# "Closure Planning and Stakeholders" row (row 8).
$ws2.Range("D8").Value = ""
$ws2.Range("E8").Value = ""
$ws2.Rows.Item(8).RowHeight = 13.8

# Give the previously-untitled Set 8 (row 9) its "Break" title.
$ws2.Range("C9").Value = "Break"

# --- Saved view state ---------------------------------------------------
# Sessions is no longer the active tab; Sets is now active with I9 selected.
[void]$ws1.Range("A21").Select()
[void]$ws2.Activate()
[void]$ws2.Range("I9").Select()
